{"js": "// Apply the tracked edit to \"Use cases.docx\":\n//  1. Remove the stray `_GoBack` bookmark that sits after the\n//     \"Software requirements:\" heading.\n//  2. Strike through four \"not implemented yet\" bullet items:\n//       - \"Get the key status/expiry date\"\n//       - \"Send a request to prolong the key\"\n//       - \"Speech-to-text support for search field\"\n//       - \"Synchronization of vocabulary and categories with the cloud on start up\"\n//  3. Re-insert the `_GoBack` bookmark around the last of those bullets\n//     (its new resting place after the edit).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Exact bullet texts (trimmed) that must be struck through.\nconst targets = [\n  \"Get the key status/expiry date\",\n  \"Send a request to prolong the key\",\n  \"Speech-to-text support for search field\",\n  \"Synchronization of vocabulary and categories with the cloud on start up\",\n];\n\n// Remove the existing _GoBack bookmark (originally right after \"Software requirements:\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nlet lastMatch = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text.trim();\n  if (targets.indexOf(text) !== -1) {\n    para.font.strikeThrough = true;\n    lastMatch = para;\n  }\n}\nawait context.sync();\n\n// Re-insert the _GoBack bookmark wrapping the text of the last struck-through\n// paragraph (\"Synchronization of vocabulary and categories with the cloud on start up\").\nif (lastMatch) {\n  lastMatch.getRange(\"Content\").insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Apply the tracked edit to \"Use cases.docx\":\n#  1. Remove the stray `_GoBack` bookmark that sits after the\n#     \"Software requirements:\" heading.\n#  2. Strike through four \"not implemented yet\" bullet items:\n#       - \"Get the key status/expiry date\"\n#       - \"Send a request to prolong the key\"\n#       - \"Speech-to-text support for search field\"\n#       - \"Synchronization of vocabulary and categories with the cloud on start up\"\n#  3. Re-insert the `_GoBack` bookmark around the last of those bullets\n#     (its new resting place after the edit).\n\n$d = $word.ActiveDocument\n\n# Remove the existing _GoBack bookmark (originally right after \"Software requirements:\").\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$targets = @(\n    \"Get the key status/expiry date\",\n    \"Send a request to prolong the key\",\n    \"Speech-to-text support for search field\",\n    \"Synchronization of vocabulary and categories with the cloud on start up\"\n)\n\n$lastRange = $null\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd()\n    if ($targets -contains $text) {\n        $p.Range.Font.StrikeThrough = $true\n        $lastRange = $p.Range\n    }\n}\n\n# Re-insert the _GoBack bookmark wrapping the last struck-through paragraph\n# (\"Synchronization of vocabulary and categories with the cloud on start up\").\nif ($lastRange -ne $null) {\n    $d.Bookmarks.Add(\"_GoBack\", $lastRange)\n}\n"}
